$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B49 currently holds the text "3" (inline string) - it should become a true
# numeric value of 3, matching the rest of the politeness_score column.
$ws.Cells.Item(49, 2).Value = 3

# Append a new annotation row (row 50) for Sunsi Wu.
$ws.Cells.Item(50, 1).Value = "Sunsi Wu"

# B50 ("politeness_score") must stay a text value "5" (not numeric), like the
# row was authored originally. Force text type via a temporary text number
# format, then restore the default style so no extra formatting is left on
# the cell.
$ws.Cells.Item(50, 2).NumberFormat = "@"
$ws.Cells.Item(50, 2).Value = "5"
$ws.Cells.Item(50, 2).Style = "Normal"

$ws.Cells.Item(50, 3).Value = "would like to thank;carefully"
$ws.Cells.Item(50, 4).Value = "SMY"
$ws.Cells.Item(50, 5).Value = "OTH"
$ws.Cells.Item(50, 6).Value = "658343d9-2c6f-4c77-9518-16756d4b8755"
$ws.Cells.Item(50, 7).Value = "SylJ1D1C-_annotated.xlsx"
$ws.Cells.Item(50, 8).Value = "First, we would like to thank the reviewer for carefully evaluating our paper."
